# The document contains two piecewise-pdf equations of the form
#   f(y) = { ... , y>0;  0, y>0. }
# where the second "0, y>0." branch is a typo and should read "y<0."
# (the complement of the first branch's condition).
#
# These "0. " / ">0. " tails live inside OMML math runs (m:t elements
# nested many levels deep in m:eqArr/m:e/m:d/m:oMath), which
# Document.Content.Find cannot reach (Find only walks the plain-text
# story, not math runs). So instead we locate every m:oMath whose raw
# OOXML contains the exact offending run, patch just that run's text
# (splitting it into "<" and "0. " to match the authored edit), and
# push the corrected fragment back with Range.InsertXML - which
# replaces the contents of the Range it's called on.
#
# We keep the m:oMathPara wrapper in the payload we feed back to
# InsertXML: the OMath.Range only covers <m:oMath>...</m:oMath>, and
# if we don't also re-supply the enclosing <m:oMathPara> the wrapper
# element gets dropped from the document on write-back.

$d = $word.ActiveDocument

# Exact OOXML of the run that needs to change: "...>0. " (with a
# leading ">" and a trailing ". ") rendered in the same run properties
# used throughout these equation arrays.
$oldRun = '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="MTSY" w:hAnsi="Cambria Math" w:cstheme="minorHAnsi"/><w:color w:val="231F20"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><m:t xml:space="preserve">&gt;0. </m:t></m:r>'

# Replacement: the same run properties, but split into a "<" run and
# a "0. " run, matching the authored diff exactly.
$newRuns = '<m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="MTSY" w:hAnsi="Cambria Math" w:cstheme="minorHAnsi"/><w:color w:val="231F20"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><m:t>&lt;</m:t></m:r><m:r><w:rPr><w:rFonts w:ascii="Cambria Math" w:eastAsia="MTSY" w:hAnsi="Cambria Math" w:cstheme="minorHAnsi"/><w:color w:val="231F20"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><m:t xml:space="preserve">0. </m:t></m:r>'

$updated = 0
$count = $d.OMaths.Count

for ($i = 1; $i -le $count; $i++) {
    $om = $d.OMaths.Item($i)
    $r = $om.Range

    # Range.WordOpenXML gives back the whole package; pull out just
    # the <m:oMath>...</m:oMath> fragment for this equation.
    $pkg = $r.WordOpenXML
    $startTag = "<m:oMath>"
    $endTag = "</m:oMath>"
    $si = $pkg.IndexOf($startTag)
    $ei = $pkg.IndexOf($endTag)
    if ($si -ge 0 -and $ei -ge 0) {
        $ei = $ei + $endTag.Length
        $fragment = $pkg.Substring($si, $ei - $si)

        if ($fragment.IndexOf($oldRun) -ge 0) {
            $newFragment = $fragment.Replace($oldRun, $newRuns)
            $payload = "<m:oMathPara>$newFragment</m:oMathPara>"
            $r.InsertXML($payload)
            $updated = $updated + 1
        }
    }
}

Write-Output "Patched $updated equation(s) with the '<0.' branch fix"
